$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column AK: header "qg21" plus values for rows 2-10
$ws.Range("AK1").Value = "qg21"

$values = @(1, 2, 3, 4, 5, 5, 4, 3, 2)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 37).Value = $values[$i]
}
